{"js": "// Auto-generated: replace 100 table-cell equation strings in row-major order.\n//\n// The document body contains one 20-row x 5-column table (100 cells, one\n// run of text per cell, e.g. \"91-76=15\"). The diff swaps each cell's text\n// for a new equation string; the Nth cell in reading order maps 1:1 onto\n// the Nth [oldText, newText] pair below (verified against the source XML).\n//\n// We replace using body.search() for the OLD text and insertText(...,\n// Word.InsertLocation.replace) on the matched Range - this only rewrites\n// the <w:t> run contents and leaves paragraph/run formatting (rFonts, sz,\n// jc) untouched, exactly matching the diff. A couple of the old equation\n// strings repeat twice in the source (e.g. \"91-76=15\"); since each\n// replacement consumes one occurrence and no replacement text ever equals\n// another old value, re-searching after every sync always lands on the\n// next remaining - i.e. earliest unreplaced - occurrence, which lines up\n// with the next cell in the pairs list.\nconst pairs = [\n  [\"91-76=15\", \"95-95=0\"],\n  [\"49-15=34\", \"54-41=13\"],\n  [\"85-74=11\", \"17+12=29\"],\n  [\"46-3=43\", \"66-7=59\"],\n  [\"13+22=35\", \"3+29=32\"],\n  [\"18+40=58\", \"83-21=62\"],\n  [\"59-4=55\", \"99-57=42\"],\n  [\"6+29=35\", \"37+54=91\"],\n  [\"47-38=9\", \"82+15=97\"],\n  [\"87-30=57\", \"41-26=15\"],\n  [\"89-67=22\", \"38+24=62\"],\n  [\"43+32=75\", \"63-44=19\"],\n  [\"6+85=91\", \"16+54=70\"],\n  [\"98-98=0\", \"38-35=3\"],\n  [\"88-54=34\", \"3+83=86\"],\n  [\"61-48=13\", \"45+45=90\"],\n  [\"71-33=38\", \"45+21=66\"],\n  [\"15+77=92\", \"73+9=82\"],\n  [\"18+20=38\", \"90-55=35\"],\n  [\"29+26=55\", \"17+17=34\"],\n  [\"83-38=45\", \"37+41=78\"],\n  [\"85-21=64\", \"95+3=98\"],\n  [\"50+42=92\", \"42+26=68\"],\n  [\"28+25=53\", \"84-11=73\"],\n  [\"33+37=70\", \"51-47=4\"],\n  [\"23-15=8\", \"53-26=27\"],\n  [\"78+8=86\", \"31+46=77\"],\n  [\"45-18=27\", \"35-2=33\"],\n  [\"55-29=26\", \"6+63=69\"],\n  [\"70+5=75\", \"38+8=46\"],\n  [\"97-28=69\", \"26+29=55\"],\n  [\"46-17=29\", \"21+13=34\"],\n  [\"24+67=91\", \"58+3=61\"],\n  [\"52+35=87\", \"93-52=41\"],\n  [\"99-27=72\", \"24+19=43\"],\n  [\"47-23=24\", \"26+48=74\"],\n  [\"6-6=0\", \"17+6=23\"],\n  [\"32-28=4\", \"75-60=15\"],\n  [\"20+15=35\", \"42-4=38\"],\n  [\"69-28=41\", \"25-0=25\"],\n  [\"6+64=70\", \"17+0=17\"],\n  [\"83-39=44\", \"31+57=88\"],\n  [\"6+13=19\", \"37+48=85\"],\n  [\"63-60=3\", \"3+2=5\"],\n  [\"64-30=34\", \"31-30=1\"],\n  [\"47-33=14\", \"95+1=96\"],\n  [\"7+52=59\", \"11+27=38\"],\n  [\"33-20=13\", \"81-78=3\"],\n  [\"26-10=16\", \"30+57=87\"],\n  [\"2+86=88\", \"20-9=11\"],\n  [\"68-30=38\", \"24-18=6\"],\n  [\"61+29=90\", \"24+57=81\"],\n  [\"27-2=25\", \"59+37=96\"],\n  [\"54+11=65\", \"50-11=39\"],\n  [\"75-48=27\", \"0+86=86\"],\n  [\"57-26=31\", \"71-46=25\"],\n  [\"14-3=11\", \"94-9=85\"],\n  [\"32-8=24\", \"35+48=83\"],\n  [\"65-60=5\", \"57-29=28\"],\n  [\"10-8=2\", \"34-26=8\"],\n  [\"16-0=16\", \"58+18=76\"],\n  [\"83-82=1\", \"43+2=45\"],\n  [\"96-92=4\", \"55-35=20\"],\n  [\"69-52=17\", \"29+51=80\"],\n  [\"45+20=65\", \"5+94=99\"],\n  [\"95-19=76\", \"26+72=98\"],\n  [\"91-64=27\", \"4+29=33\"],\n  [\"5+56=61\", \"15+76=91\"],\n  [\"87+5=92\", \"21+54=75\"],\n  [\"85-72=13\", \"23+42=65\"],\n  [\"54-32=22\", \"28-11=17\"],\n  [\"26+61=87\", \"44+48=92\"],\n  [\"7+63=70\", \"99-38=61\"],\n  [\"49-20=29\", \"64+15=79\"],\n  [\"14+82=96\", \"22+61=83\"],\n  [\"2+0=2\", \"73+4=77\"],\n  [\"91-76=15\", \"49-21=28\"],\n  [\"39-20=19\", \"22+63=85\"],\n  [\"59-58=1\", \"93-37=56\"],\n  [\"91-52=39\", \"23-21=2\"],\n  [\"9+59=68\", \"67-1=66\"],\n  [\"89-57=32\", \"11+35=46\"],\n  [\"35+44=79\", \"2+47=49\"],\n  [\"82+4=86\", \"81-43=38\"],\n  [\"20+12=32\", \"67-22=45\"],\n  [\"96-28=68\", \"94-55=39\"],\n  [\"45-30=15\", \"80-51=29\"],\n  [\"68+10=78\", \"99-55=44\"],\n  [\"17+20=37\", \"7+12=19\"],\n  [\"21-19=2\", \"24+6=30\"],\n  [\"94-89=5\", \"55-15=40\"],\n  [\"54-34=20\", \"24+69=93\"],\n  [\"11+36=47\", \"98-51=47\"],\n  [\"83-35=48\", \"18+11=29\"],\n  [\"33+0=33\", \"2+10=12\"],\n  [\"68-30=38\", \"65-41=24\"],\n  [\"4+5=9\", \"90-57=33\"],\n  [\"93-90=3\", \"39+8=47\"],\n  [\"89-86=3\", \"29+37=66\"],\n  [\"8+26=34\", \"55+28=83\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length < 1) {\n    throw new Error(`Occurrence of \"${oldText}\" not found`);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Auto-generated: replace 100 table-cell equation strings in row-major order.\n# The table has 20 rows x 5 columns = 100 cells; the Nth cell (reading order)\n# maps 1:1 to the Nth (oldText,newText) pair below (verified against the source\n# XML). We replace using Find.Execute(..., Replace:=1) (wdReplaceOne) on a fresh\n# range anchored at the document start each time, so duplicate strings (e.g.\n# \"91-76=15\" appears twice) are replaced one-by-one, earliest-remaining-occurrence\n# first - matching the cell order in the table. MatchCase/MatchWholeWord are on so\n# no text is a partial match of another; the run/paragraph formatting (rFonts, sz,\n# jc) is preserved untouched, only the <w:t> content changes, exactly like the diff.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('91-76=15', '95-95=0'),\n    @('49-15=34', '54-41=13'),\n    @('85-74=11', '17+12=29'),\n    @('46-3=43', '66-7=59'),\n    @('13+22=35', '3+29=32'),\n    @('18+40=58', '83-21=62'),\n    @('59-4=55', '99-57=42'),\n    @('6+29=35', '37+54=91'),\n    @('47-38=9', '82+15=97'),\n    @('87-30=57', '41-26=15'),\n    @('89-67=22', '38+24=62'),\n    @('43+32=75', '63-44=19'),\n    @('6+85=91', '16+54=70'),\n    @('98-98=0', '38-35=3'),\n    @('88-54=34', '3+83=86'),\n    @('61-48=13', '45+45=90'),\n    @('71-33=38', '45+21=66'),\n    @('15+77=92', '73+9=82'),\n    @('18+20=38', '90-55=35'),\n    @('29+26=55', '17+17=34'),\n    @('83-38=45', '37+41=78'),\n    @('85-21=64', '95+3=98'),\n    @('50+42=92', '42+26=68'),\n    @('28+25=53', '84-11=73'),\n    @('33+37=70', '51-47=4'),\n    @('23-15=8', '53-26=27'),\n    @('78+8=86', '31+46=77'),\n    @('45-18=27', '35-2=33'),\n    @('55-29=26', '6+63=69'),\n    @('70+5=75', '38+8=46'),\n    @('97-28=69', '26+29=55'),\n    @('46-17=29', '21+13=34'),\n    @('24+67=91', '58+3=61'),\n    @('52+35=87', '93-52=41'),\n    @('99-27=72', '24+19=43'),\n    @('47-23=24', '26+48=74'),\n    @('6-6=0', '17+6=23'),\n    @('32-28=4', '75-60=15'),\n    @('20+15=35', '42-4=38'),\n    @('69-28=41', '25-0=25'),\n    @('6+64=70', '17+0=17'),\n    @('83-39=44', '31+57=88'),\n    @('6+13=19', '37+48=85'),\n    @('63-60=3', '3+2=5'),\n    @('64-30=34', '31-30=1'),\n    @('47-33=14', '95+1=96'),\n    @('7+52=59', '11+27=38'),\n    @('33-20=13', '81-78=3'),\n    @('26-10=16', '30+57=87'),\n    @('2+86=88', '20-9=11'),\n    @('68-30=38', '24-18=6'),\n    @('61+29=90', '24+57=81'),\n    @('27-2=25', '59+37=96'),\n    @('54+11=65', '50-11=39'),\n    @('75-48=27', '0+86=86'),\n    @('57-26=31', '71-46=25'),\n    @('14-3=11', '94-9=85'),\n    @('32-8=24', '35+48=83'),\n    @('65-60=5', '57-29=28'),\n    @('10-8=2', '34-26=8'),\n    @('16-0=16', '58+18=76'),\n    @('83-82=1', '43+2=45'),\n    @('96-92=4', '55-35=20'),\n    @('69-52=17', '29+51=80'),\n    @('45+20=65', '5+94=99'),\n    @('95-19=76', '26+72=98'),\n    @('91-64=27', '4+29=33'),\n    @('5+56=61', '15+76=91'),\n    @('87+5=92', '21+54=75'),\n    @('85-72=13', '23+42=65'),\n    @('54-32=22', '28-11=17'),\n    @('26+61=87', '44+48=92'),\n    @('7+63=70', '99-38=61'),\n    @('49-20=29', '64+15=79'),\n    @('14+82=96', '22+61=83'),\n    @('2+0=2', '73+4=77'),\n    @('91-76=15', '49-21=28'),\n    @('39-20=19', '22+63=85'),\n    @('59-58=1', '93-37=56'),\n    @('91-52=39', '23-21=2'),\n    @('9+59=68', '67-1=66'),\n    @('89-57=32', '11+35=46'),\n    @('35+44=79', '2+47=49'),\n    @('82+4=86', '81-43=38'),\n    @('20+12=32', '67-22=45'),\n    @('96-28=68', '94-55=39'),\n    @('45-30=15', '80-51=29'),\n    @('68+10=78', '99-55=44'),\n    @('17+20=37', '7+12=19'),\n    @('21-19=2', '24+6=30'),\n    @('94-89=5', '55-15=40'),\n    @('54-34=20', '24+69=93'),\n    @('11+36=47', '98-51=47'),\n    @('83-35=48', '18+11=29'),\n    @('33+0=33', '2+10=12'),\n    @('68-30=38', '65-41=24'),\n    @('4+5=9', '90-57=33'),\n    @('93-90=3', '39+8=47'),\n    @('89-86=3', '29+37=66'),\n    @('8+26=34', '55+28=83')\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 1)\n    if (-not $found) {\n        throw \"Occurrence of $oldText not found\"\n    }\n}\n\n"}
